# Weekly update: insert a new "Apio" price record (Primera + Segunda quality
# rows) ahead of the existing row 323, pushing the rest of the log down by
# two rows (A1:R412 -> A1:R414).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the position of the current row 323; Excel will
# carry the date-format style of row 322 down onto the freshly inserted
# row 323/324 (same behaviour as manually inserting rows in the UI).
$ws.Rows("323:324").Insert()

# --- New row 323 (Calidad = Primera) ---
$ws.Cells.Item(323, 1).Value2  = 8
$ws.Cells.Item(323, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(323, 3).Value2  = "Coquimbo"
$ws.Cells.Item(323, 4).Value2  = 44711
$ws.Cells.Item(323, 5).Value2  = 4
$ws.Cells.Item(323, 6).Value2  = 100112017
$ws.Cells.Item(323, 7).Value2  = "Apio"
$ws.Cells.Item(323, 8).Value2  = "Americana (o)"
$ws.Cells.Item(323, 9).Value2  = "Primera"
$ws.Cells.Item(323, 10).Value2 = 2500
$ws.Cells.Item(323, 11).Value2 = 8000
$ws.Cells.Item(323, 12).Value2 = 9000
$ws.Cells.Item(323, 13).Value2 = 8500
$ws.Cells.Item(323, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(323, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(323, 16).Value2 = 1417
$ws.Cells.Item(323, 17).Value2 = 6
$ws.Cells.Item(323, 18).Value2 = "Hortaliza"

# --- New row 324 (Calidad = Segunda) ---
$ws.Cells.Item(324, 1).Value2  = 8
$ws.Cells.Item(324, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(324, 3).Value2  = "Coquimbo"
$ws.Cells.Item(324, 4).Value2  = 44711
$ws.Cells.Item(324, 5).Value2  = 4
$ws.Cells.Item(324, 6).Value2  = 100112017
$ws.Cells.Item(324, 7).Value2  = "Apio"
$ws.Cells.Item(324, 8).Value2  = "Americana (o)"
$ws.Cells.Item(324, 9).Value2  = "Segunda"
$ws.Cells.Item(324, 10).Value2 = 1500
$ws.Cells.Item(324, 11).Value2 = 6000
$ws.Cells.Item(324, 12).Value2 = 7000
$ws.Cells.Item(324, 13).Value2 = 6500
$ws.Cells.Item(324, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(324, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(324, 16).Value2 = 1083
$ws.Cells.Item(324, 17).Value2 = 6
$ws.Cells.Item(324, 18).Value2 = "Hortaliza"
